$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "2024-12-31 21:00:00"
$ws.Cells.Item(4, 2).Value = "2025-01-01 02:30:00"
$ws.Cells.Item(4, 3).Value = 250.42
$ws.Cells.Item(4, 4).Value = 20784.86
$ws.Cells.Item(4, 5).Value = -147.74
$ws.Cells.Item(4, 6).Value = -0.71
$ws.Cells.Item(4, 7).Value = 21022.24
$ws.Cells.Item(4, 8).Value = 20702.69
$ws.Cells.Item(4, 9).Value = 20952.52
$ws.Cells.Item(4, 10).Value = 20932.6
